$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2052
$ws1.Range("F7").Value = 3375
$ws1.Range("F9").Value = 823

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 2

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2
$ws4.Range("F6").Value = 2052
$ws4.Range("F8").Value = 3375
$ws4.Range("F10").Value = 823
